$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Insert a new row above row 15 (shifts existing rows 15-35 down to 16-36),
# adjusting formulas, shared formula ranges, merged cells, conditional
# formatting ranges and the used range automatically.
$ws.Rows.Item(15).Insert()

# The freshly inserted row doesn't pick up the surrounding cell formatting,
# so copy it down from the row below (which holds the formerly-row-15
# "Circuit diagrams/Calculations" task, now shifted to row 16).
$ws.Range("A16:BK16").Copy()
$ws.Range("A15:BK15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(15).RowHeight = 30

# Fill in the new task row: "Cardboard Prototype"
$ws.Cells.Item(15, 2).Value = "Cardboard Prototype"
$ws.Cells.Item(15, 4).Value = 44088
$ws.Cells.Item(15, 5).Formula = "=D15+5"

# Update the frozen-pane scroll position and active selection.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B14").Select()
